$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Column L = Seasonality Index updates
$ws1.Range("L2").Value = 0.8100000000000001
$ws1.Range("L3").Value = 1.17
$ws1.Range("L4").Value = 1.02
$ws1.Range("L5").Value = 0.82

# Row 6: MyForecast, Inventory Coverage, Stockout Risk, Seasonality Index
$ws1.Range("D6").Value = 33
$ws1.Range("H6").Value = 0.54
$ws1.Range("I6").Value = "Low"
$ws1.Range("L6").Value = 1.17

# Row 7
$ws1.Range("D7").Value = 29
$ws1.Range("L7").Value = 1

# Row 8
$ws1.Range("D8").Value = 32
$ws1.Range("L8").Value = 1

# Row 9
$ws1.Range("D9").Value = 33
$ws1.Range("L9").Value = 1.05

# Row 10
$ws1.Range("D10").Value = 34
$ws1.Range("L10").Value = 0.98

# Row 11
$ws1.Range("D11").Value = 33
$ws1.Range("L11").Value = 0.88

# Row 12
$ws1.Range("D12").Value = 32
$ws1.Range("L12").Value = 0.86

# Row 13
$ws1.Range("D13").Value = 33
$ws1.Range("L13").Value = 0.84

# Row 14
$ws1.Range("L14").Value = 1.02

# Row 15
$ws1.Range("L15").Value = 0.96

# Row 16
$ws1.Range("L16").Value = 1.04

# Row 17
$ws1.Range("L17").Value = 0.88

# --- Sheet: Summary ---
# These cells hold their numbers as text (inline strings) in the workbook,
# so force text format before assigning to avoid Excel auto-converting to a number.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "535"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "271"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "36"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "29"
